# Fix the "Evidenciado" level labels in column B of the
# "Resp Responde pelo cumprimento " sheet: the original labels were
# written without a space ("MuitoEvidenciado", "PoucoEvidenciado",
# "NãoEvidenciado") while the rest of the workbook (column A) already
# used the spaced form. Re-write column B so both columns agree.
#
# Row 11's "Não Evidenciado" additionally picked up a trailing space in
# the source data, so it is special-cased to keep that exact text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Resp Responde pelo cumprimento ")

$ws.Range("B3").Value  = "Muito Evidenciado"
$ws.Range("B5").Value  = "Pouco Evidenciado"
$ws.Range("B6").Value  = "Não Evidenciado"

$ws.Range("B8").Value  = "Muito Evidenciado"
$ws.Range("B10").Value = "Pouco Evidenciado"
$ws.Range("B11").Value = "Não Evidenciado "

$ws.Range("B13").Value = "Muito Evidenciado"
$ws.Range("B15").Value = "Pouco Evidenciado"
$ws.Range("B16").Value = "Não Evidenciado"

$ws.Range("B18").Value = "Muito Evidenciado"
$ws.Range("B20").Value = "Pouco Evidenciado"
$ws.Range("B21").Value = "Não Evidenciado"

$ws.Range("B23").Value = "Muito Evidenciado"
$ws.Range("B25").Value = "Pouco Evidenciado"
$ws.Range("B26").Value = "Não Evidenciado"

# Match the saved view: sheet 2 active/selected with B26 as the current
# selection and the view scrolled back to the top-left corner.
$ws.Activate()
$ws.Range("B26").Select()
